# iAMC_Kdr_kinetics.xlsx edit
#
# Summary of the change being applied (per the target diff):
#  - Each of the three sheets (Kdr_amp, Kdr_tau_act, Kdr_tau_deact) has a
#    per-sheet row label in A1 ("Vmem (mV)") and a per-sheet unique label in
#    A2 ("Kdr_amplitude (pA/pF)", "Kdr_tau_activation (ms)",
#    "Kdr_tau_deactivation (ms)"). The edit blanks out A1 and renames every
#    A2 label to the single shared word "mean" (collapsing the shared
#    string table down to one string).
#  - The special bold/underlined/red font used for the A2 label cells loses
#    its bold, underline and red-color formatting (becomes a plain font).
#  - Sheet selection/active-cell bookkeeping changes: Kdr_amp's sheet view
#    no longer has "tabSelected" and its selection collapses from the whole
#    used range down to the next data cell (A2); Kdr_tau_act's and
#    Kdr_tau_deact's selections likewise collapse to A3; and the workbook's
#    active tab moves from Kdr_amp (index 0) to Kdr_tau_deact (index 2),
#    which becomes the tab marked "tabSelected".
#  - Kdr_amp additionally gains an explicit portrait/A4 page setup.

$wb = $excel.ActiveWorkbook

$wsAmp   = $wb.Worksheets.Item(1)   # Kdr_amp
$wsAct   = $wb.Worksheets.Item(2)   # Kdr_tau_act
$wsDeact = $wb.Worksheets.Item(3)   # Kdr_tau_deact

# --- Relabel the header/summary cells on every sheet -----------------------
$sheets = @($wsAmp, $wsAct, $wsDeact)
foreach ($ws in $sheets) {
    $ws.Range("A1").Value = ""
    $ws.Range("A2").Value = "mean"
}

# --- Strip bold / underline / red color from the A2 label font -------------
# (done as three separate passes over all three sheets so the identical
# resulting font/style is shared by every sheet instead of being re-created
# per sheet)
foreach ($ws in $sheets) {
    $ws.Range("A2").Font.Bold = $false
}
foreach ($ws in $sheets) {
    $ws.Range("A2").Font.Underline = $false
}
foreach ($ws in $sheets) {
    $ws.Range("A2").Font.ThemeColor = 1
}

# --- Kdr_amp gains an explicit page setup (portrait, A4) -------------------
$wsAmp.PageSetup.PaperSize = 9
$wsAmp.PageSetup.Orientation = 1

# --- Update per-sheet selections --------------------------------------------
$wsAmp.Activate()
$wsAmp.Range("A2").Select()

$wsAct.Activate()
$wsAct.Range("A3").Select()

# Kdr_tau_deact ends up both selected-within and the active workbook tab.
$wsDeact.Activate()
$wsDeact.Range("A3").Select()
